# Update "paises.xlsx" (Pais sheet) with the refreshed COVID-19 scrape:
#  - a few countries changed rank order (new si text per fixed row), and
#  - case/recovered/death counters were refreshed for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Column A text / country name updates (ranking reshuffle) ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 13 de Junio de 2020 a las 22:25'
$ws.Cells.Item(31, 1).Value = 'Egipto'
$ws.Cells.Item(32, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(78, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(79, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(80, 1).Value = 'Guinea'
$ws.Cells.Item(81, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(144, 1).Value = 'Ruanda'
$ws.Cells.Item(145, 1).Value = 'Malaui'
$ws.Cells.Item(146, 1).Value = 'Togo'
$ws.Cells.Item(206, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(207, 1).Value = 'Groenlandia'
$ws.Cells.Item(210, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 1).Value = 'Seychelles'
$ws.Cells.Item(213, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(214, 1).Value = 'Papua Nueva Guinea'

# --- Numeric updates (refreshed case counts) ---
$ws.Cells.Item(4, 2).Value = 2136198
$ws.Cells.Item(4, 3).Value = 19276
$ws.Cells.Item(4, 4).Value = 847064
$ws.Cells.Item(4, 5).Value = 1171748
$ws.Cells.Item(4, 7).Value = 561
$ws.Cells.Item(4, 8).Value = 117386
$ws.Cells.Item(5, 2).Value = 832866
$ws.Cells.Item(5, 3).Value = 2964
$ws.Cells.Item(5, 5).Value = 363201
$ws.Cells.Item(5, 7).Value = 154
$ws.Cells.Item(5, 8).Value = 42055
$ws.Cells.Item(7, 2).Value = 321626
$ws.Cells.Item(7, 3).Value = 12023
$ws.Cells.Item(7, 5).Value = 150095
$ws.Cells.Item(12, 2).Value = 187420
$ws.Cells.Item(12, 3).Value = 169
$ws.Cells.Item(12, 5).Value = 6653
$ws.Cells.Item(12, 7).Value = 4
$ws.Cells.Item(12, 8).Value = 8867
$ws.Cells.Item(24, 2).Value = 65736
$ws.Cells.Item(24, 3).Value = 3809
$ws.Cells.Item(24, 4).Value = 36850
$ws.Cells.Item(24, 5).Value = 27463
$ws.Cells.Item(24, 7).Value = 69
$ws.Cells.Item(24, 8).Value = 1423
$ws.Cells.Item(30, 2).Value = 46356
$ws.Cells.Item(30, 3).Value = 578
$ws.Cells.Item(30, 4).Value = 22865
$ws.Cells.Item(30, 5).Value = 19617
$ws.Cells.Item(30, 7).Value = 46
$ws.Cells.Item(30, 8).Value = 3874
$ws.Cells.Item(31, 2).Value = 42980
$ws.Cells.Item(31, 3).Value = 1677
$ws.Cells.Item(31, 4).Value = 11529
$ws.Cells.Item(31, 5).Value = 29967
$ws.Cells.Item(31, 7).Value = 62
$ws.Cells.Item(31, 8).Value = 1484
$ws.Cells.Item(32, 2).Value = 41990
$ws.Cells.Item(32, 3).Value = 491
$ws.Cells.Item(32, 4).Value = 26761
$ws.Cells.Item(32, 5).Value = 14941
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 288
$ws.Cells.Item(73, 2).Value = 7007
$ws.Cells.Item(73, 3).Value = 128
$ws.Cells.Item(73, 4).Value = 2556
$ws.Cells.Item(73, 5).Value = 4004
$ws.Cells.Item(73, 7).Value = 14
$ws.Cells.Item(73, 8).Value = 447
$ws.Cells.Item(78, 2).Value = 4848
$ws.Cells.Item(78, 3).Value = 164
$ws.Cells.Item(78, 4).Value = 2397
$ws.Cells.Item(78, 5).Value = 2406
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 45
$ws.Cells.Item(79, 2).Value = 4724
$ws.Cells.Item(79, 3).Value = 87
$ws.Cells.Item(79, 4).Value = 595
$ws.Cells.Item(79, 5).Value = 4023
$ws.Cells.Item(79, 7).Value = 5
$ws.Cells.Item(79, 8).Value = 106
$ws.Cells.Item(80, 2).Value = 4484
$ws.Cells.Item(80, 3).Value = 58
$ws.Cells.Item(80, 4).Value = 3213
$ws.Cells.Item(80, 5).Value = 1246
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = 25
$ws.Cells.Item(81, 2).Value = 4449
$ws.Cells.Item(81, 3).Value = 8
$ws.Cells.Item(81, 4).Value = 2823
$ws.Cells.Item(81, 5).Value = 1585
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 41
$ws.Cells.Item(131, 5).Value = 135
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 14
$ws.Cells.Item(140, 2).Value = 659
$ws.Cells.Item(140, 3).Value = 20
$ws.Cells.Item(140, 4).Value = 176
$ws.Cells.Item(140, 5).Value = 471
$ws.Cells.Item(144, 2).Value = 541
$ws.Cells.Item(144, 3).Value = 31
$ws.Cells.Item(144, 4).Value = 332
$ws.Cells.Item(144, 5).Value = 207
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 2
$ws.Cells.Item(145, 2).Value = 529
$ws.Cells.Item(145, 3).Value = 48
$ws.Cells.Item(145, 4).Value = 66
$ws.Cells.Item(145, 5).Value = 458
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 5
$ws.Cells.Item(146, 2).Value = 525
$ws.Cells.Item(146, 4).Value = 279
$ws.Cells.Item(146, 5).Value = 233
$ws.Cells.Item(146, 8).Value = 13
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0

Write-Host "Updated countries & provincias Spain data (13 Jun 2020, 22:25 refresh)."
